$d = $word.ActiveDocument

# --- 1. Remove the existing "_GoBack" bookmark (previously sitting right
#        after "...manera correcta" and before the following comma). ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2. Locate the bold "Humedad" heading so we can split the run that
#        currently holds "med" into three runs: "m", "e", "d". ---
$findRange = $d.Content.Duplicate
$found = $findRange.Find.Execute("Humedad:", $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0)

$wordStart = $findRange.Start
# "Humedad:" -> H(0) u(1) m(2) e(3) d(4) a(5) d(6) :(7)
$mStart = $wordStart + 2
$mEnd   = $wordStart + 3
$eStart = $wordStart + 3
$eEnd   = $wordStart + 4

# Toggling Bold off/on (net no-op) forces Word to break the run boundary
# at these single-character ranges, splitting "med" into "m" | "e" | "d".
$rngM = $d.Range($mStart, $mEnd)
$rngM.Bold = 0
$rngM.Bold = 1

$rngE = $d.Range($eStart, $eEnd)
$rngE.Bold = 0
$rngE.Bold = 1

# --- 3. Re-insert the "_GoBack" bookmark between the new "e" and "d" runs. ---
$insertPoint = $d.Range($eEnd, $eEnd)
$d.Bookmarks.Add("_GoBack", $insertPoint)
